$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates: "Girassol" -> "Sunflower" in the two rows that were
# already adjusted (these keep their original cell style, s="2") ---
$ws.Cells.Item(24, 3).Value = 'Ela foi no Sunflower todo mês\ncom as cartas e presentes da minha mãe.'
$ws.Cells.Item(27, 3).Value = 'Você é do Sunflower?'

# --- Style fix-up: column C translation cells should use the same style
# as column B (s="1") instead of the mismatched s="2" style, for every
# translated row except the two just-adjusted ones (24 and 27) which
# intentionally retain their original style. ---
$rows = @(2,3,4,5,7,8,9,10,11,12,13,14,15,17,18,19,20,22,23,25,29,30,31,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
